$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 72, shifting existing rows 72:181 down to 73:182.
$ws.Rows.Item(72).Insert()

# Populate the new row 72 with data (matches the constant pattern used throughout
# this table for columns A, B, C, E, F, G, H, I, N, Q, R, plus the new values for
# D, J, K, L, M, O, P).
$ws.Range("A72").Value = 10
$ws.Range("B72").Value = "Vega Modelo de Temuco"
$ws.Range("C72").Value = "La Araucanía"
$ws.Range("D72").Value = 44797
$ws.Range("E72").Value = 9
$ws.Range("F72").Value = 100112012
$ws.Range("G72").Value = "Espinaca"
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 30
$ws.Range("K72").Value = 9000
$ws.Range("L72").Value = 9000
$ws.Range("M72").Value = 9000
$ws.Range("N72").Value = "`$/docena de atados"
$ws.Range("O72").Value = "Región de Coquimbo"
$ws.Range("P72").Value = 3000
$ws.Range("Q72").Value = 3
$ws.Range("R72").Value = "Hortaliza"

# Ensure the date cell keeps the same date/time number format as the rest of
# column D.
$ws.Range("D72").NumberFormat = "YYYY-MM-DD HH:MM:SS"
